$wb = $excel.ActiveWorkbook

# Sheet "G1" (second sheet) gets a "download" value of 1 for rows 2-46 and 49-55
$ws = $wb.Worksheets.Item("G1")
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 2).Value = 1
}
for ($r = 49; $r -le 55; $r++) {
    $ws.Cells.Item($r, 2).Value = 1
}

# Selection on G1 sheet and scroll position
$ws.Range("C6").Select()

# Switch active/selected tab to "G2" (third sheet)
$ws2 = $wb.Worksheets.Item("G2")
$ws2.Activate()
$ws2.Range("C21").Select()
